$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 446, shifting existing rows 446:503 down to 447:504
$ws.Rows.Item(446).Insert()

# Populate the new row 446 with this week's data (same market/region/category as neighboring rows)
$ws.Range("A446").Value2 = 4
$ws.Range("B446").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C446").Value2 = "Los Lagos"
$ws.Range("D446").Value2 = 45131
$ws.Range("D446").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E446").Value2 = 10
$ws.Range("F446").Value2 = 100112017
$ws.Range("G446").Value2 = "Apio"
$ws.Range("H446").Value2 = "Americana (o)"
$ws.Range("I446").Value2 = "Primera"
$ws.Range("J446").Value2 = 20
$ws.Range("K446").Value2 = 11000
$ws.Range("L446").Value2 = 11000
$ws.Range("M446").Value2 = 11000
$ws.Range("N446").Value2 = "$/docena de matas"
$ws.Range("O446").Value2 = "Región de Coquimbo"
$ws.Range("P446").Value2 = 1833
$ws.Range("Q446").Value2 = 6
$ws.Range("R446").Value2 = "Hortaliza"
